$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets (trailing "Every 250/500/1000 ..." titles were cleaned up)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Every 250 Operating Hours (~Mon").Name = "Every 250 OP Hrs(Monthly)"
$wb.Worksheets.Item("Every 500 Hours (~3 Months)").Name = "Every 500 Hours (3 Months)"
$wb.Worksheets.Item("Every 1000 Hours   Yearly").Name = "Every 1000 Hours (Yearly)"

# ---------------------------------------------------------------------------
# 2. "Daily" sheet - append a new "DRUM MIXER" block (rows 20-27)
# ---------------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily")
$dailyNew = @(
    "DRUM MIXER :",
    "Check drum for concrete build-up: Clean any leftover mix after each use",
    "Inspect charging and discharging hoppers: Ensure no blockages or cracks",
    "Check drum rotation (smooth and consistent): Listen for abnormal noises",
    "Verify water tank level: Sufficient for drum washing",
    "Check chute condition and movement: No cracks or deformations",
    "Look for hydraulic leaks near drum motor and lines: Report if found",
    "Ensure drum cover is in place: Prevents concrete splashing during transport"
)
$row = 20
foreach ($txt in $dailyNew) {
    $wsDaily.Cells.Item(2, 1).Copy()
    $wsDaily.Cells.Item($row, 1).PasteSpecial(-4122)
    $wsDaily.Cells.Item($row, 1).Value = $txt
    $row++
}

# ---------------------------------------------------------------------------
# 3. "Weekly" sheet - row 11 text swap + new rows 12-17
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly")
$wsWeekly.Cells.Item(11, 1).Value = "DRUM MIXER:"
$weeklyNew = @(
    "Inspect drum blades (mixing fins): Look for wear, cracks, or loosening",
    "Check drum shell for dents or cracks: Especially near welds and seams",
    "Inspect drum ring gear & drive pinion: Clean and visually check for damage or misalignment",
    "Check drum lock mechanism (if present): Should engage securely",
    "Inspect chute lifting mechanism: Hydraulic or mechanical—check cables, hinges, cylinders",
    " Lubricate mixer according to Maintenance Checklist in this section."
)
$row = 12
foreach ($txt in $weeklyNew) {
    $wsWeekly.Cells.Item(5, 1).Copy()
    $wsWeekly.Cells.Item($row, 1).PasteSpecial(-4122)
    $wsWeekly.Cells.Item($row, 1).Value = $txt
    $row++
}

# ---------------------------------------------------------------------------
# 4. "Every 250 OP Hrs(Monthly)" sheet - new rows 14-19
# ---------------------------------------------------------------------------
$ws250 = $wb.Worksheets.Item("Every 250 OP Hrs(Monthly)")
$new250 = @(
    "DRUM MIXER:",
    "Check drum wall thickness (if accessible): Tap test or ultrasonic check (in high-wear fleets)",
    "Inspect drum rotation speed and drive motor performance: Confirm proper RPM and torque under load",
    "Examine hydraulic oil level and cleanliness: Top up and change if dark or foamy",
    "Check water spray nozzles: Ensure they’re unclogged and spraying properly",
    "Inspect chute pivots and discharge door: Lubricate and check for misalignment or jamming"
)
$row = 14
foreach ($txt in $new250) {
    $ws250.Cells.Item(6, 1).Copy()
    $ws250.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws250.Cells.Item($row, 1).Value = $txt
    $row++
}

# ---------------------------------------------------------------------------
# 5. "Every 500 Hours (3 Months)" sheet
#    row7 / row8 text swap (gain the shared style) + new row9 (no style) +
#    new rows 10-13 (styled)
# ---------------------------------------------------------------------------
$ws500 = $wb.Worksheets.Item("Every 500 Hours (3 Months)")

$ws500.Cells.Item(6, 1).Copy()
$ws500.Cells.Item(7, 1).PasteSpecial(-4122)
$ws500.Cells.Item(7, 1).Value = "DRUM MIXER:"

$ws500.Cells.Item(6, 1).Copy()
$ws500.Cells.Item(8, 1).PasteSpecial(-4122)
$ws500.Cells.Item(8, 1).Value = "Grease rear manual controller."

$ws500.Cells.Item(1, 1).Copy()
$ws500.Cells.Item(9, 1).PasteSpecial(-4122)
$ws500.Cells.Item(9, 1).Value = " Change hydraulic oil for temperature reasons."

$new500 = @(
    "Perform detailed inspection of blade wear: Measure remaining blade height; replace if below limit",
    "Inspect drum trunnion rollers (support rollers):  Check for flat spots or excessive play",
    "Inspect and tighten drive coupling bolts: Drum drive to hydraulic motor",
    "Lubricate gear teeth and roller bearings: Use appropriate high-temp grease"
)
$row = 10
foreach ($txt in $new500) {
    $ws500.Cells.Item(6, 1).Copy()
    $ws500.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws500.Cells.Item($row, 1).Value = $txt
    $row++
}

# ---------------------------------------------------------------------------
# 6. "Every 1000 Hours (Yearly)" sheet - new rows 12-17
# ---------------------------------------------------------------------------
$ws1000 = $wb.Worksheets.Item("Every 1000 Hours (Yearly)")
$new1000 = @(
    "DRUM MIXER:",
    "Perform non-destructive testing (NDT) on drum shell: Ultrasonic or visual if in high-hour units",
    "Replace or rebuild mixing blades if needed: Based on wear percentage or volume mixed",
    "Flush and replace hydraulic oil and filters: Ensure clean system operation",
    "Inspect entire water system (tank, pump, piping): Flush and sanitize if needed",
    "Repaint or treat drum exterior for corrosion: Especially important for coastal or humid areas"
)
$row = 12
foreach ($txt in $new1000) {
    $ws1000.Cells.Item(6, 1).Copy()
    $ws1000.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws1000.Cells.Item($row, 1).Value = $txt
    $row++
}
